# Regenerate s_vals data to filter save games.
# Updates columns B, C, D, E, G for rows 2-6 on the active sheet with
# the newly regenerated values (G = B + C + D + E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 3.272327238179451;  C = 1.626987699542094; D = 3.223369029078222;  E = 0.5333859586016987; G = 8.656069925401464 }
    3 = @{ B = 3.272327238179451;  C = 1.626987699542094; D = 0.7210945179870265; E = 0.5333859586016987; G = 6.15379541431027 }
    4 = @{ B = 3.272327238179451;  C = 1.626987699542094; D = 0.1496068669990043; E = 0.5333859586016987; G = 5.582307763322248 }
    5 = @{ B = 3.272327238179451;  C = 1.626987699542094; D = 0.7210945179870265; E = 0.5333859586016987; G = 6.15379541431027 }
    6 = @{ B = 3.272327238179451;  C = 1.626987699542094; D = 0.7210945179870265; E = 0.5333859586016987; G = 6.15379541431027 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
